$wb = $excel.ActiveWorkbook

# --- Sheet "001": new (smaller-magnitude) uncertainty template values ---
$ws1 = $wb.Worksheets.Item("001")
$ws1.Range("D2:M2").Value = -0.05
$ws1.Range("N2:Y2").Value = -0.01
$ws1.Range("D3:M3").Value = -0.025
$ws1.Range("N3:Y3").Value = -0.005
$ws1.Range("L9").Select()

# --- Sheet "002": new uncertainty template values (positive side already 0) ---
$ws2 = $wb.Worksheets.Item("002")
$ws2.Range("N2:Y2").Value = -0.05
$ws2.Range("N3:Y3").Value = -0.01
$ws2.Range("D4:M4").Value = 0.005
$ws2.Range("N4:Y4").Value = 0.01
$ws2.Range("D5:M5").Value = 0.01
$ws2.Range("N5:Y5").Value = 0.05
$ws2.Range("G10").Select()

# --- Sheet "003": new uncertainty template values ---
$ws3 = $wb.Worksheets.Item("003")
$ws3.Range("N2:Y2").Value = -0.01
$ws3.Range("N3:Y3").Value = -0.005
$ws3.Range("D4:M4").Value = 0.01
$ws3.Range("N4:Y4").Value = 0.005
$ws3.Range("D5:M5").Value = 0.025
$ws3.Range("N5:Y5").Value = 0.01
$ws3.Range("Z7").Select()

# --- Sheet "004": initial data already all 0 stocks, only the selection moves ---
$ws4 = $wb.Worksheets.Item("004")
$ws4.Range("F15").Select()

# "003" is the sheet that was active/saved last in the source workbook.
$ws3.Activate()
